$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per row, reflecting re-run of NATMI pipeline (Dr Hou advice)
# Columns: E=Ligand-expressing cells, G/H=Ligand avg/total expr, I/J=Ligand derived specificity avg/total,
#          K=Receptor-expressing cells, M/N=Receptor avg/total expr, O/P=Receptor derived specificity avg/total,
#          Q/R=Edge avg/total expr weight, S/T=Edge avg/total expr derived specificity

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 7.727457666666666
$ws.Cells.Item(2, 8).Value = 23.182373
$ws.Cells.Item(2, 9).Value = 0.1630271452636819
$ws.Cells.Item(2, 10).Value = 0.1630271452636819
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 36.802266
$ws.Cells.Item(2, 14).Value = 110.406798
$ws.Cells.Item(2, 15).Value = 0.4381152826760633
$ws.Cells.Item(2, 16).Value = 0.4381152826760634
$ws.Cells.Item(2, 17).Value = 284.3879525524059
$ws.Cells.Item(2, 18).Value = 2559.491572971654
$ws.Cells.Item(2, 19).Value = 0.07142468383106965
$ws.Cells.Item(2, 20).Value = 0.07142468383106965

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 7.727457666666666
$ws.Cells.Item(3, 8).Value = 23.182373
$ws.Cells.Item(3, 9).Value = 0.1630271452636819
$ws.Cells.Item(3, 10).Value = 0.1630271452636819
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 34.02833166666667
$ws.Cells.Item(3, 14).Value = 102.084995
$ws.Cells.Item(3, 15).Value = 0.4050927773614947
$ws.Cells.Item(3, 16).Value = 0.4050927773614947
$ws.Cells.Item(3, 17).Value = 262.9524924214594
$ws.Cells.Item(3, 18).Value = 2366.572431793134
$ws.Cells.Item(3, 19).Value = 0.06604111906018076
$ws.Cells.Item(3, 20).Value = 0.06604111906018074

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 7.727457666666666
$ws.Cells.Item(4, 8).Value = 23.182373
$ws.Cells.Item(4, 9).Value = 0.1630271452636819
$ws.Cells.Item(4, 10).Value = 0.1630271452636819
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 13.17073133333333
$ws.Cells.Item(4, 14).Value = 39.512194
$ws.Cells.Item(4, 15).Value = 0.156791939962442
$ws.Cells.Item(4, 16).Value = 0.156791939962442
$ws.Cells.Item(4, 17).Value = 101.7762688173735
$ws.Cells.Item(4, 18).Value = 915.9864193563619
$ws.Cells.Item(4, 19).Value = 0.02556134237243152
$ws.Cells.Item(4, 20).Value = 0.02556134237243152

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 37.42779933333333
$ws.Cells.Item(5, 8).Value = 112.283398
$ws.Cells.Item(5, 9).Value = 0.7896189849264272
$ws.Cells.Item(5, 10).Value = 0.7896189849264271
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 36.802266
$ws.Cells.Item(5, 14).Value = 110.406798
$ws.Cells.Item(5, 15).Value = 0.4381152826760633
$ws.Cells.Item(5, 16).Value = 0.4381152826760634
$ws.Cells.Item(5, 17).Value = 1377.427826859956
$ws.Cells.Item(5, 18).Value = 12396.8504417396
$ws.Cells.Item(5, 19).Value = 0.3459441447874279
$ws.Cells.Item(5, 20).Value = 0.3459441447874279

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 37.42779933333333
$ws.Cells.Item(6, 8).Value = 112.283398
$ws.Cells.Item(6, 9).Value = 0.7896189849264272
$ws.Cells.Item(6, 10).Value = 0.7896189849264271
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 34.02833166666667
$ws.Cells.Item(6, 14).Value = 102.084995
$ws.Cells.Item(6, 15).Value = 0.4050927773614947
$ws.Cells.Item(6, 16).Value = 0.4050927773614947
$ws.Cells.Item(6, 17).Value = 1273.605569268112
$ws.Cells.Item(6, 18).Value = 11462.45012341301
$ws.Cells.Item(6, 19).Value = 0.3198689476612106
$ws.Cells.Item(6, 20).Value = 0.3198689476612105

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 37.42779933333333
$ws.Cells.Item(7, 8).Value = 112.283398
$ws.Cells.Item(7, 9).Value = 0.7896189849264272
$ws.Cells.Item(7, 10).Value = 0.7896189849264271
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 13.17073133333333
$ws.Cells.Item(7, 14).Value = 39.512194
$ws.Cells.Item(7, 15).Value = 0.156791939962442
$ws.Cells.Item(7, 16).Value = 0.156791939962442
$ws.Cells.Item(7, 17).Value = 492.9514894172458
$ws.Cells.Item(7, 18).Value = 4436.563404755212
$ws.Cells.Item(7, 19).Value = 0.1238058924777887
$ws.Cells.Item(7, 20).Value = 0.1238058924777887

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 2.244565
$ws.Cells.Item(8, 8).Value = 6.733695
$ws.Cells.Item(8, 9).Value = 0.04735386980989085
$ws.Cells.Item(8, 10).Value = 0.04735386980989083
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 36.802266
$ws.Cells.Item(8, 14).Value = 110.406798
$ws.Cells.Item(8, 15).Value = 0.4381152826760633
$ws.Cells.Item(8, 16).Value = 0.4381152826760634
$ws.Cells.Item(8, 17).Value = 82.60507818429
$ws.Cells.Item(8, 18).Value = 743.44570365861
$ws.Cells.Item(8, 19).Value = 0.02074645405756583
$ws.Cells.Item(8, 20).Value = 0.02074645405756583

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 2.244565
$ws.Cells.Item(9, 8).Value = 6.733695
$ws.Cells.Item(9, 9).Value = 0.04735386980989085
$ws.Cells.Item(9, 10).Value = 0.04735386980989083
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 34.02833166666667
$ws.Cells.Item(9, 14).Value = 102.084995
$ws.Cells.Item(9, 15).Value = 0.4050927773614947
$ws.Cells.Item(9, 16).Value = 0.4050927773614947
$ws.Cells.Item(9, 17).Value = 76.37880226739168
$ws.Cells.Item(9, 18).Value = 687.4092204065249
$ws.Cells.Item(9, 19).Value = 0.01918271064010332
$ws.Cells.Item(9, 20).Value = 0.01918271064010331

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 2.244565
$ws.Cells.Item(10, 8).Value = 6.733695
$ws.Cells.Item(10, 9).Value = 0.04735386980989085
$ws.Cells.Item(10, 10).Value = 0.04735386980989083
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 13.17073133333333
$ws.Cells.Item(10, 14).Value = 39.512194
$ws.Cells.Item(10, 15).Value = 0.156791939962442
$ws.Cells.Item(10, 16).Value = 0.156791939962442
$ws.Cells.Item(10, 17).Value = 29.56256257520334
$ws.Cells.Item(10, 18).Value = 266.06306317683
$ws.Cells.Item(10, 19).Value = 0.007424705112221698
$ws.Cells.Item(10, 20).Value = 0.007424705112221697
